# Auto-generated edit script applying the symbol-list refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that hold plain (non numeric-looking) text - assign directly.
$textUpdates = @{
    'E18' = '17OneONEWorstin24h'
    'B41' = 'KickToken'
    'C41' = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
    'E41' = '40KickTokenKICK'
    'B43' = 'BKEXToken'
    'C43' = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
    'E43' = '42BKEXTokenBKK'
}
foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

# Cells that hold numeric-looking text (prices in column D, hour counters in
# column G). The source workbook stores these as text, so force the text
# number format before writing the new value - otherwise Excel's normal
# auto-detection would silently convert them to real numbers.
$numericTextUpdates = @{
    'D2' = '243.45'
    'G2' = '8'
    'D3' = '23.72'
    'G3' = '8'
    'D4' = '5.290'
    'G4' = '8'
    'D5' = '0.05792'
    'G5' = '8'
    'D6' = '6.480'
    'G6' = '8'
    'D7' = '3.329'
    'G7' = '8'
    'D8' = '0.8102'
    'G8' = '8'
    'D9' = '0.8736'
    'G9' = '8'
    'D10' = '0.1383'
    'G10' = '8'
    'G11' = '8'
    'D12' = '0.03084'
    'G12' = '8'
    'D13' = '0.03056'
    'G13' = '8'
    'D14' = '0.09323'
    'G14' = '8'
    'D15' = '3.865'
    'G15' = '8'
    'D16' = '0.001534'
    'G16' = '8'
    'D17' = '0.04695'
    'G17' = '8'
    'D18' = '0.0006019'
    'G18' = '8'
    'D19' = '0.006119'
    'G19' = '8'
    'D20' = '0.001300'
    'G20' = '8'
    'D21' = '0.004598'
    'G21' = '8'
    'D22' = '0.00008701'
    'G22' = '8'
    'G23' = '8'
    'D24' = '2.142'
    'G24' = '8'
    'D25' = '0.3212'
    'G25' = '8'
    'G26' = '8'
    'G27' = '8'
    'D28' = '0.0002343'
    'G28' = '8'
    'G29' = '8'
    'G30' = '8'
    'G31' = '8'
    'G32' = '8'
    'G33' = '8'
    'G34' = '8'
    'G35' = '8'
    'G36' = '8'
    'G37' = '8'
    'G38' = '8'
    'G39' = '8'
    'D40' = '0.03775'
    'G40' = '8'
    'D41' = '0.006371'
    'G41' = '8'
    'D42' = '0.004000'
    'G42' = '8'
    'D43' = '0.1053'
    'G43' = '8'
    'D44' = '0.007644'
    'G44' = '8'
    'D45' = '0.00005476'
    'G45' = '8'
    'G46' = '8'
    'D47' = '0.5899'
    'G47' = '8'
    'D48' = '0.004275'
    'G48' = '8'
    'D49' = '0.00002100'
    'G49' = '8'
    'D50' = '0.0002000'
    'G50' = '8'
    'G51' = '8'
}
foreach ($addr in $numericTextUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextUpdates[$addr]
}
